# Chen_2002: Updated USDA soil type
#
# Adds a new controlled-vocabulary field "pro_usda_soil_order" to the
# workbook:
#   - "profile" sheet: insert a new column N (pro_usda_soil_order) before
#     the existing pro_soil_taxon column, and populate the USDA soil
#     order ("Oxisols") + WRB equivalent ("Ferrasols") + taxonomy system
#     ("WRB") for the three data rows.
#   - "controlled vocabulary" sheet: insert a new column E before the
#     existing pro_soil_taxon_sys column, with header
#     "pro_usda_soil_order" and the 12 USDA soil order values as the
#     allowed list.

$wb = $excel.ActiveWorkbook

$profile = $wb.Worksheets.Item("profile")
$profile.Columns.Item(14).Insert()

$cv = $wb.Worksheets.Item("controlled vocabulary")
$cv.Columns.Item(5).Insert()

# ---------------------------------------------------------------------
# 1) "profile" sheet — header + first WRB->USDA mapping value
# ---------------------------------------------------------------------
$profile.Range("O4").Value = "Ferrasols"
$profile.Range("N1").Value = "pro_usda_soil_order"

# ---------------------------------------------------------------------
# 2) "controlled vocabulary" sheet — new field name + allowed values
# ---------------------------------------------------------------------
$cv.Range("E2").Value = "pro_usda_soil_order"

$cv.Range("E4").Value = "Alfisols"
$cv.Range("E5").Value = "Andisols"
$cv.Range("E6").Value = "Aridisols"
$cv.Range("E7").Value = "Entisols"
$cv.Range("E8").Value = "Gelisols"
$cv.Range("E9").Value = "Histosols"
$cv.Range("E10").Value = "Inceptisols"
$cv.Range("E11").Value = "Mollisols"
$cv.Range("E12").Value = "Oxisols"
$cv.Range("E13").Value = "Spodosols"
$cv.Range("E14").Value = "Ultisols"
$cv.Range("E15").Value = "Vertisols"

# ---------------------------------------------------------------------
# 3) "profile" sheet — remaining data rows (reuse existing strings)
# ---------------------------------------------------------------------
$profile.Range("N4").Value = "Oxisols"
$profile.Range("Q4").Value = "WRB"

$profile.Range("N5").Value = "Oxisols"
$profile.Range("O5").Value = "Ferrasols"
$profile.Range("Q5").Value = "WRB"

$profile.Range("N6").Value = "Oxisols"
$profile.Range("O6").Value = "Ferrasols"
$profile.Range("Q6").Value = "WRB"
